# Reorder the "Recorded By" list (column G) on every session row so that
# any entry equal to "System" (case-insensitive) that appears alongside
# other recorders is moved out of the front position, effectively
# reversing the order of the whole comma-separated list.
#
# Example: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#          "system, System, backup@backdoor.com" -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ",\s*"

    if ($parts.Count -le 1) {
        continue
    }

    $hasSystemToken = $false
    foreach ($p in $parts) {
        if ($p.Trim().ToLower() -eq "system") {
            $hasSystemToken = $true
        }
    }

    if ($hasSystemToken) {
        $reversedParts = $parts[($parts.Count - 1)..0]
        $newVal = [string]::Join(", ", $reversedParts)
        $cell.Value2 = $newVal
    }
}
